# Completed 6 user stories
# Expand the "Acceptance Tests" table (3rd table in the document) with the
# full set of acceptance-test bullet points for each user story, and add
# four brand new user-story rows (rebate process, add a pet, submit
# receipt, continue through form).

$d = $word.ActiveDocument
$t = $d.Tables.Item(3)

# --- Row 2: Sarah / "Sarah would like to log in." ------------------------
$row2 = $t.Rows.Item(2)
$row2.Cells.Item(3).Range.Text = "If Sarah enters in an incorrect email, verify that an error messages appears telling her that the email is not recognized. `rIf Sarah enters in an incorrect password, verify that an error message appears.`rIf Sarah enters in valid details, verify that she is logged in."

# --- Row 3: Alex / "Alex would like to create an account." ---------------
$row3 = $t.Rows.Item(3)
$row3.Cells.Item(3).Range.Text = "If Alex enters in an invalid address, verify an error message appears.`rIf Alex enters in an invalid email, verify an error message appears.`rIf Alex enters an invalid phone number, verify an error message appears.`rIf Alex’s details are correct, verify he is logged in"

# --- New row: Sarah / Alex would like to start the rebate process. -------
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Sarah"
$newRow.Cells.Item(2).Range.Text = "Alex would like to start the rebate process."
$newRow.Cells.Item(3).Range.Text = "If the date is invalid, verify that the button is disabled.`rIf the offer code is invalid, verify that an error message appears.`rIf the offer code and date are valid, verify that Alex is taken to a specific offer page"

# --- New row: Alex / Alex would like to add a pet. ------------------------
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Alex"
$newRow.Cells.Item(2).Range.Text = "Alex would like to add a pet."
$newRow.Cells.Item(3).Range.Text = "If the pet’s birth date is invalid, verify that the button is disabled.`rIf the pet’s animal is invalid, verify an error message is displayed.`rIf the user leaves a form empty, verify an error message is displayed.`rIf all the details are valid, verify a pet is added to their account."

# --- New row: Sarah / Sarah would like to submit her receipt -------------
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Sarah "
$newRow.Cells.Item(2).Range.Text = "Sarah would like to submit her receipt"
$newRow.Cells.Item(3).Range.Text = "If Alex submits an invalid image, verify that it is not accepted.`rIf Alex submits a valid image and the clinic is detected correctly, verify that she is taken to the next step.`rIf Alex submits a valid images but the clinic is not detected, verify that she will be able to change it."

# --- New row: Sarah / Sarah would like to continue through the form ------
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Sarah "
$newRow.Cells.Item(2).Range.Text = "Sarah would like to continue through the form"
$newRow.Cells.Item(3).Range.Text = "If Sarah has redeemed a rebate before verify that her card is chosen by default and can be changed`rIf Sarah has used a pet before, verify that this pet will be the default."
